$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.810.83'
$ws.Range('E2').Value = '  +6.48%  '

$ws.Range('D3').Value = '3.487.02'
$ws.Range('E3').Value = '  +7.42%  '

$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').Value = '585.72'
$ws.Range('E5').Value = '  +7.55%  '

$ws.Range('D6').Value = '159.86'
$ws.Range('E6').Value = '  +8.30%  '

$ws.Range('E7').Value = '  -0.17%  '

$ws.Range('D8').Value = '3.490.28'
$ws.Range('E8').Value = '  +7.43%  '

$ws.Range('E9').Value = '  +2.78%  '

$ws.Range('E10').Value = '  +2.67%  '

$ws.Range('E11').Value = '  +8.42%  '

$ws.Range('E12').Value = '  +3.05%  '

$ws.Range('D13').Value = '4.087.66'
$ws.Range('E13').Value = '  +7.15%  '

$ws.Range('E14').Value = '  -0.24%  '

$ws.Range('E15').Value = '  +9.77%  '

$ws.Range('D16').Value = '27.83'
$ws.Range('E16').Value = '  +5.48%  '

$ws.Range('D17').Value = '64.825.18'
$ws.Range('E17').Value = '  +6.35%  '

$ws.Range('D18').Value = '3.480.66'
$ws.Range('E18').Value = '  +6.73%  '

$ws.Range('E19').Value = '  +2.66%  '

$ws.Range('D20').Value = '14.41'
$ws.Range('E20').Value = '  +7.36%  '

$ws.Range('D21').Value = '400.68'
$ws.Range('E21').Value = '  +6.31%  '

$ws.Range('D22').Value = '8.60'
$ws.Range('E22').Value = '  +2.20%  '

$ws.Range('E23').Value = '  +3.48%  '

$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.20%  '

$ws.Range('E25').Value = '  +3.42%  '

$ws.Range('E26').Value = '  +23.63%  '

$ws.Range('D27').Value = '9.71'
$ws.Range('E27').Value = '  +12.70%  '

$ws.Range('E28').Value = '  +6.53%  '

$ws.Range('D29').Value = '1.01'
$ws.Range('E29').Value = '  +0.55%  '

$ws.Range('E30').Value = '  +13.75%  '

$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').Value = '6.74'
$ws.Range('E31').Value = '  +8.99%  '

$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = '5.93'
$ws.Range('E32').Value = '  +9.76%  '

$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '2.07'
$ws.Range('E33').Value = '  +6.91%  '

$ws.Range('E34').Value = '  +5.93%  '

$ws.Range('E35').Value = '  -0.07%  '

$ws.Range('D36').Value = '7.00'
$ws.Range('E36').Value = '  +5.46%  '

$ws.Range('D37').Value = '1.52'
$ws.Range('E37').Value = '  +5.15%  '

$ws.Range('D38').Value = '158.18'
$ws.Range('E38').Value = '  -0.91%  '

$ws.Range('D39').Value = '28.67'
$ws.Range('E39').Value = '  +8.89%  '

$ws.Range('E40').Value = '  +10.81%  '

$ws.Range('D41').Value = '0.0782'
$ws.Range('E41').Value = '  +8.52%  '

$ws.Range('D42').Value = '2.922.32'
$ws.Range('E42').Value = '  +4.97%  '

$ws.Range('E43').Value = '  +3.66%  '

$ws.Range('E44').Value = '  +7.33%  '

$ws.Range('E45').Value = '  +4.26%  '

$ws.Range('D46').Value = '41.93'
$ws.Range('E46').Value = '  +4.63%  '

$ws.Range('E47').Value = '  +10.50%  '

$ws.Range('E48').Value = '  +6.94%  '

$ws.Range('D49').Value = '3.535.25'
$ws.Range('E49').Value = '  +7.19%  '

$ws.Range('E50').Value = '  +25.63%  '

$ws.Range('D51').Value = '6.51'
$ws.Range('E51').Value = '  +4.48%  '
